$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header labels in row 1: "_old" suffix -> "_FV2410", "_new" suffix -> "_FV2504"
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value2 = ($cell.Value2 -replace "_old$", "_FV2410")
}

foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value2 = ($cell.Value2 -replace "_new$", "_FV2504")
}

# 2) Convert the data range into an Excel Table (ListObject) with autofilter
$dataRange = $ws.Range("A1:U73")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"

# 3) Freeze the top (header) row
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
